$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icam2"
$ws.Range("C2").Value = "Itgal"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 27.94460233333334
$ws.Range("H2").Value = 83.83380700000001
$ws.Range("I2").Value = 0.9231270698117087
$ws.Range("J2").Value = 0.9231270698117087
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 21.087087
$ws.Range("N2").Value = 63.261261
$ws.Range("O2").Value = 0.980881822159902
$ws.Range("P2").Value = 0.980881822159902
$ws.Range("Q2").Value = 589.2702605834031
$ws.Range("R2").Value = 5303.432345250627
$ws.Range("S2").Value = 0.9054785623220399
$ws.Range("T2").Value = 0.9054785623220399

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icam2"
$ws.Range("C3").Value = "Itgal"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 27.94460233333334
$ws.Range("H3").Value = 83.83380700000001
$ws.Range("I3").Value = 0.9231270698117087
$ws.Range("J3").Value = 0.9231270698117087
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3230143333333333
$ws.Range("N3").Value = 0.9690430000000001
$ws.Range("O3").Value = 0.01502525635066456
$ws.Range("P3").Value = 0.01502525635066456
$ws.Range("Q3").Value = 9.02650709296678
$ws.Range("R3").Value = 81.23856383670102
$ws.Range("S3").Value = 0.01387022086815874
$ws.Range("T3").Value = 0.01387022086815874

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Icam2"
$ws.Range("C4").Value = "Itgal"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 27.94460233333334
$ws.Range("H4").Value = 83.83380700000001
$ws.Range("I4").Value = 0.9231270698117087
$ws.Range("J4").Value = 0.9231270698117087
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.08799
$ws.Range("N4").Value = 0.26397
$ws.Range("O4").Value = 0.00409292148943331
$ws.Range("P4").Value = 0.004092921489433309
$ws.Range("Q4").Value = 2.45884555931
$ws.Range("R4").Value = 22.12961003379
$ws.Range("S4").Value = 0.003778286621509946
$ws.Range("T4").Value = 0.003778286621509945

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icam2"
$ws.Range("C5").Value = "Itgal"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.091026
$ws.Range("H5").Value = 3.273078
$ws.Range("I5").Value = 0.03604115107650029
$ws.Range("J5").Value = 0.03604115107650029
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 21.087087
$ws.Range("N5").Value = 63.261261
$ws.Range("O5").Value = 0.980881822159902
$ws.Range("P5").Value = 0.980881822159902
$ws.Range("Q5").Value = 23.006560181262
$ws.Range("R5").Value = 207.059041631358
$ws.Range("S5").Value = 0.03535210994065792
$ws.Range("T5").Value = 0.03535210994065792

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Icam2"
$ws.Range("C6").Value = "Itgal"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.091026
$ws.Range("H6").Value = 3.273078
$ws.Range("I6").Value = 0.03604115107650029
$ws.Range("J6").Value = 0.03604115107650029
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.3230143333333333
$ws.Range("N6").Value = 0.9690430000000001
$ws.Range("O6").Value = 0.01502525635066456
$ws.Range("P6").Value = 0.01502525635066456
$ws.Range("Q6").Value = 0.3524170360393334
$ws.Range("R6").Value = 3.171753324354
$ws.Range("S6").Value = 0.0005415275340974468
$ws.Range("T6").Value = 0.0005415275340974468

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Icam2"
$ws.Range("C7").Value = "Itgal"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.091026
$ws.Range("H7").Value = 3.273078
$ws.Range("I7").Value = 0.03604115107650029
$ws.Range("J7").Value = 0.03604115107650029
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.08799
$ws.Range("N7").Value = 0.26397
$ws.Range("O7").Value = 0.00409292148943331
$ws.Range("P7").Value = 0.004092921489433309
$ws.Range("Q7").Value = 0.09599937774
$ws.Range("R7").Value = 0.8639943996599999
$ws.Range("S7").Value = 0.0001475136017449205
$ws.Range("T7").Value = 0.0001475136017449205

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Icam2"
$ws.Range("C8").Value = "Itgal"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.236046333333333
$ws.Range("H8").Value = 3.708139
$ws.Range("I8").Value = 0.04083177911179101
$ws.Range("J8").Value = 0.04083177911179101
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 21.087087
$ws.Range("N8").Value = 63.261261
$ws.Range("O8").Value = 0.980881822159902
$ws.Range("P8").Value = 0.980881822159902
$ws.Range("Q8").Value = 26.064616567031
$ws.Range("R8").Value = 234.581549103279
$ws.Range("S8").Value = 0.04005114989720419
$ws.Range("T8").Value = 0.04005114989720419

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Icam2"
$ws.Range("C9").Value = "Itgal"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.236046333333333
$ws.Range("H9").Value = 3.708139
$ws.Range("I9").Value = 0.04083177911179101
$ws.Range("J9").Value = 0.04083177911179101
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3230143333333333
$ws.Range("N9").Value = 0.9690430000000001
$ws.Range("O9").Value = 0.01502525635066456
$ws.Range("P9").Value = 0.01502525635066456
$ws.Range("Q9").Value = 0.3992606823307778
$ws.Range("R9").Value = 3.593346140977
$ws.Range("S9").Value = 0.0006135079484083703
$ws.Range("T9").Value = 0.0006135079484083703

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Icam2"
$ws.Range("C10").Value = "Itgal"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.236046333333333
$ws.Range("H10").Value = 3.708139
$ws.Range("I10").Value = 0.04083177911179101
$ws.Range("J10").Value = 0.04083177911179101
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.08799
$ws.Range("N10").Value = 0.26397
$ws.Range("O10").Value = 0.00409292148943331
$ws.Range("P10").Value = 0.004092921489433309
$ws.Range("Q10").Value = 0.10875971687
$ws.Range("R10").Value = 0.9788374518299999
$ws.Range("S10").Value = 0.0001671212661784436
$ws.Range("T10").Value = 0.0001671212661784435
